$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "[1] First, it’s important to understand how your teen’s brain is changing. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[1] Kwekucala, kubalulekile kucondza kutsi ingcondvo yemntfwana wakho igucuka njani. ",
    2)

$d.Content.Find.Execute(
    "Your teen’s brain hasn’t changed this fast since they were a baby. It is normal for teens to behave in risky ways as they push the limits of what is possible and allowed, especially if it is pleasurable or fun. This is normal behaviour. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kutfutfukisa Ingcondvo yemntfwana loseminyakeni yekutfomba ayikaze igucuke ngekushesha lokunjalo kusukela aseluswane. Kuvamile kutsi bantfu batiphatse ngendlela lengaba yingoti njengobe cindzetela imincele yaloko lokungenteka nalokuvunyelwe, ikakhulukati nangabe tijabulisa. Loku kutiphatsa kuvamile. ",
    2)

$d.Content.Find.Execute(
    "[3] To support your teen through these changes, you can guide their choices by balancing rules and limits with independence. Try giving them more freedom to make their own choices when they show they can be responsible and take responsibility for their actions. Most importantly: be supportive, loving, and positive with your teen. This helps them feel safe and secure around you, and they're more likely to come to you for help or advice when they need it.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[3] Kute usite umntfwana wakho loseminyakeni yekutfomba kutsi akhone kwenta lolushintjo, ungamcondzisa endleleni lokufanele akhetse ngayo ngekutsi ulinganisele emkhatsini wemitsetfo kanye nemincele kanye nekukhululeka. Yetama kubaniketa inkhululeko yekwenta tincumo tabo nangabe bakhombisa kutsi bangakhona kutiphatsa futsi batitsatsele imitfwalo yetento tabo. Lokubaluleke kakhulu kutsi usekele, utsandze futsi ube nembono lomuhle ngemntfwana wakho. Loku kuyabasita kutsi bative baphephile nabakuwe, futsi batawuta kuwe kute batfole lusito nobe seluleko nabalidinga.",
    2)
